# Electricity sector updates:
#  1) Update offshore wind capacity to align with latest EU targets
#  2) Modify electricity parameters (onshore wind, solar PV) to better
#     calibrate the power sector
#
# Sheet "BGDPbES" rows:
#   7  = onshore wind
#   8  = solar PV
#   15 = offshore wind
# Columns B:AK are the year columns (2015-2050).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BGDPbES")

# Make this the active / selected sheet (mirrors the author switching to it)
$ws.Activate()

# Set guaranteed dispatch percentage to 1 (100%) across all years
$ws.Range("B7:AK7").Value = 1
$ws.Range("B8:AK8").Value = 1
$ws.Range("B15:AK15").Value = 1

# Leave the selection where the author apparently ended up working
$ws.Range("J36").Select()
